# Append updated COVID data rows (465-491) to Sheet1, matching the
# "aggiornamento fino a 6 gennaio 2022" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(465, 44539, 4, 10, 267.6659528907923),
    @(466, 44540, 1, 9, 240.8993576017131),
    @(467, 44541, 1, 9, 240.8993576017131),
    @(468, 44542, 4, 13, 347.9657387580299),
    @(469, 44543, 0, 11, 294.4325481798715),
    @(470, 44544, 2, 12, 321.1991434689508),
    @(471, 44545, 0, 12, 321.1991434689508),
    @(472, 44546, 0, 8, 214.1327623126338),
    @(473, 44547, 1, 8, 214.1327623126338),
    @(474, 44548, 5, 12, 321.1991434689508),
    @(475, 44550, 1, 9, 240.8993576017131),
    @(476, 44551, 6, 15, 401.4989293361884),
    @(477, 44552, 2, 15, 401.4989293361884),
    @(478, 44553, 2, 17, 455.0321199143469),
    @(479, 44554, 1, 18, 481.7987152034261),
    @(480, 44555, 5, 22, 588.865096359743),
    @(481, 44556, 7, 24, 642.3982869379015),
    @(482, 44557, 5, 28, 749.4646680942184),
    @(483, 44558, 1, 23, 615.6316916488222),
    @(484, 44559, 2, 23, 615.6316916488222),
    @(485, 44560, 7, 28, 749.4646680942184),
    @(486, 44561, 11, 38, 1017.130620985011),
    @(487, 44562, 9, 42, 1124.197002141328),
    @(488, 44563, 0, 35, 936.8308351177729),
    @(489, 44564, 3, 33, 883.2976445396147),
    @(490, 44565, 3, 35, 936.8308351177729),
    @(491, 44566, 2, 35, 936.8308351177729)
)

# The last pre-existing row (464) carries the A-column date style (border +
# centered + YYYY-MM-DD HH:MM:SS format) that needs to be replicated onto
# every new row's A cell, same as Excel's "fill down" / copy-format would do.
$lastRow = 464

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($lastRow, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = 0
